# Refresh the Molecule list in column A with the new source data
# (new/updated entries replacing the old sample, extended from 5 to 9 rows),
# matching the "Queue size for K8s" source-data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newValues = @(
    "Atracurium Besylate 25 Mg /2.5Ml",
    "Avibactam-0.5Gm + Ceftazidime-2Gm",
    "Azithromycin (500Mg)",
    "Aztreonam (1Gm)",
    "Bacillus Clausii 2 Billion",
    "Bandage",
    "Bed &Pillow Cover",
    "Bed Bath",
    "Benfotiamine 150 Mg+Elemental Iron 100 Mg+Elemental Zinc 7.5 Mg+L-Methylfolate 1 Mg+Methylcobalamin 1500 Mcg+Pyridoxine 10 Mg+Vitamin D3 1000 Iu"
)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $newValues[$i]
}

# The old rows (2-6) carried the highlighted "data row" style; the refreshed
# rows use the worksheet's default (unstyled) formatting, same as the newly
# appended rows (7-10).
$ws.Range("A2:A10").Style = "Normal"

# Update the selection to match where the edit left off.
$ws.Range("A6").Select() | Out-Null

Write-Output "Refreshed Molecule list: $($newValues.Length) data rows (A2:A$(1 + $newValues.Length))."
